$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update classification result values (rows 4-7, columns B-M)
# Row 4 (linear)
$ws.Range("B4").Value = 0.6248218701785994
$ws.Range("C4").Value = 0.575
$ws.Range("D4").Value = 0.7027318861503392
$ws.Range("E4").Value = 0.6615
$ws.Range("F4").Value = 0.5684898867295347
$ws.Range("G4").Value = 0.5679999999999999
$ws.Range("H4").Value = 0.5733349751257749
$ws.Range("I4").Value = 0.571
$ws.Range("J4").Value = 0.9432157734716838
$ws.Range("K4").Value = 0.9410000000000001
$ws.Range("L4").Value = 0.9494730136612695
$ws.Range("M4").Value = 0.9429999999999999

# Row 5 (rbf)
$ws.Range("B5").Value = 0.6419547449854721
$ws.Range("C5").Value = 0.623
$ws.Range("D5").Value = 0.6730249821976331
$ws.Range("E5").Value = 0.6575
$ws.Range("F5").Value = 0.503783930320159
$ws.Range("G5").Value = 0.4799999999999999
$ws.Range("H5").Value = 0.5718238318515813
$ws.Range("I5").Value = 0.5405
$ws.Range("J5").Value = 0.848534569367612
$ws.Range("K5").Value = 0.849
$ws.Range("L5").Value = 0.8678155414158395
$ws.Range("M5").Value = 0.8505

# Row 6 (poly)
$ws.Range("B6").Value = 0.4870089134085142
$ws.Range("C6").Value = 0.3810000000000001
$ws.Range("D6").Value = 0.7168958043509163
$ws.Range("E6").Value = 0.611
$ws.Range("F6").Value = 0.5624108185560841
$ws.Range("G6").Value = 0.5599999999999999
$ws.Range("H6").Value = 0.5695412193233521
$ws.Range("I6").Value = 0.5679999999999999
$ws.Range("J6").Value = 0.8103591405528215
$ws.Range("K6").Value = 0.8219999999999998
$ws.Range("L6").Value = 0.8214003094457821
$ws.Range("M6").Value = 0.8099999999999999

# Row 7 (sigmoid)
$ws.Range("B7").Value = 0.5511428465469944
$ws.Range("C7").Value = 0.5570000000000001
$ws.Range("D7").Value = 0.5557461496970595
$ws.Range("E7").Value = 0.5615000000000001
$ws.Range("F7").Value = 0.6668896321070233
$ws.Range("G7").Value = 1
$ws.Range("H7").Value = 0.500251256281407
$ws.Range("I7").Value = 0.5004999999999999
$ws.Range("J7").Value = 0.8720395474396281
$ws.Range("K7").Value = 0.8699999999999999
$ws.Range("L7").Value = 0.8875172878895017
$ws.Range("M7").Value = 0.8739999999999999
